$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 408, shifting existing rows 408-444 down to 409-445
$ws.Rows.Item(408).Insert()

# Populate the newly inserted row 408 with the new weekly record
$ws.Cells.Item(408, 1).Value = 7
$ws.Cells.Item(408, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(408, 3).Value = "Ñuble"
$ws.Cells.Item(408, 4).Value = 45166
$ws.Cells.Item(408, 5).Value = 16
$ws.Cells.Item(408, 6).Value = 100112017
$ws.Cells.Item(408, 7).Value = "Apio"
$ws.Cells.Item(408, 8).Value = "Americana (o)"
$ws.Cells.Item(408, 9).Value = "Primera"
$ws.Cells.Item(408, 10).Value = 280
$ws.Cells.Item(408, 11).Value = 6000
$ws.Cells.Item(408, 12).Value = 7000
$ws.Cells.Item(408, 13).Value = 6357
$ws.Cells.Item(408, 14).Value = "$/docena de matas"
$ws.Cells.Item(408, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(408, 16).Value = 1060
$ws.Cells.Item(408, 17).Value = 6
$ws.Cells.Item(408, 18).Value = "Hortaliza"
